$wb = $excel.ActiveWorkbook

# The edited data lives on "Sheet3" (the 3rd worksheet / active tab)
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

# FirstName: Rony -> Anthony
$ws.Range("B2").Value = "Anthony"

# NationalID: 9514931229018 -> 9414931229018
$ws.Range("F2").Value = "9414931229018"

# SmartCardNo: 9200963821 -> 9200963827
$ws.Range("G2").Value = "9200963827"

# PassportNo: B00229480 -> B00229880
$ws.Range("H2").Value = "B00229880"

# Update the active cell selection on Sheet3 from H5 to H2
$ws.Range("H2").Select()
